# Insert one new data row at row 228 (pushing the existing row 228..325
# down to 229..326, carrying all of their original values with them),
# then populate the new row 228 with the new weekly record.
#
# Columns A,B,C,E,F,G,H,I,J,Q,R,T are constant across every data row in
# this sheet, so the new row reuses those same constant values; only
# D (Fecha), M (Volumen), N (Precio minimo), O (Precio maximo),
# P (Precio promedio ponderado) and S (Precio $/Kg) carry genuinely new
# data, while K (Variedad) and L (Calidad) keep the values that already
# occupied row 228 before the insert ("Start Ruby" / "Primera").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(228).Insert()

$ws.Range("A228").Value = 4
$ws.Range("B228").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C228").Value = "Los Lagos"
$ws.Range("D228").Value = 44755
$ws.Range("E228").Value = 10
$ws.Range("F228").Value = "Fruta"
$ws.Range("G228").Value = 100102
$ws.Range("H228").Value = "Cítricos"
$ws.Range("I228").Value = 100102006
$ws.Range("J228").Value = "Pomelo"
$ws.Range("K228").Value = "Start Ruby"
$ws.Range("L228").Value = "Primera"
$ws.Range("M228").Value = 80
$ws.Range("N228").Value = 14000
$ws.Range("O228").Value = 14000
$ws.Range("P228").Value = 14000
$ws.Range("Q228").Value = "$/caja 14 kilos empedrada"
$ws.Range("R228").Value = "Región de O'Higgins"
$ws.Range("S228").Value = 1000
$ws.Range("T228").Value = 14
